$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45188 -> 45189) for every data row (rows 2 through 481).
$ws.Range("C2:C481").Value = 45189
